$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency prices and 1h volume percentages.
# Rows 38 and 39 also swap their Coin name / Link (InternetComputer(DFINITY) <-> FraxShare).
$ws.Cells.Item(2, 4).Value = '28.584.91'
$ws.Cells.Item(2, 5).Value = '  +4.42%  '
$ws.Cells.Item(3, 4).Value = '1.792.19'
$ws.Cells.Item(3, 5).Value = '  +0.79%  '
$c = $ws.Cells.Item(4, 4); $c.NumberFormat = "@"; $c.Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  +0.21%  '
$c = $ws.Cells.Item(5, 4); $c.NumberFormat = "@"; $c.Value = '313.92'
$ws.Cells.Item(5, 5).Value = '  +0.07%  '
$c = $ws.Cells.Item(6, 4); $c.NumberFormat = "@"; $c.Value = '1.003'
$ws.Cells.Item(6, 5).Value = '  +0.30%  '
$c = $ws.Cells.Item(7, 4); $c.NumberFormat = "@"; $c.Value = '0.5347'
$ws.Cells.Item(7, 5).Value = '  -0.35%  '
$c = $ws.Cells.Item(8, 4); $c.NumberFormat = "@"; $c.Value = '0.3808'
$ws.Cells.Item(8, 5).Value = '  +0.77%  '
$c = $ws.Cells.Item(9, 4); $c.NumberFormat = "@"; $c.Value = '0.07516'
$ws.Cells.Item(9, 5).Value = '  +1.40%  '
$c = $ws.Cells.Item(10, 4); $c.NumberFormat = "@"; $c.Value = '42.50'
$ws.Cells.Item(10, 5).Value = '  -0.84%  '
$ws.Cells.Item(11, 5).Value = '  +2.12%  '
$c = $ws.Cells.Item(12, 4); $c.NumberFormat = "@"; $c.Value = '1.005'
$ws.Cells.Item(12, 5).Value = '  +0.46%  '
$c = $ws.Cells.Item(13, 4); $c.NumberFormat = "@"; $c.Value = '21.07'
$ws.Cells.Item(13, 5).Value = '  +1.95%  '
$c = $ws.Cells.Item(14, 4); $c.NumberFormat = "@"; $c.Value = '6.181'
$ws.Cells.Item(14, 5).Value = '  +1.33%  '
$c = $ws.Cells.Item(15, 4); $c.NumberFormat = "@"; $c.Value = '7.395'
$ws.Cells.Item(15, 5).Value = '  +5.81%  '
$ws.Cells.Item(16, 4).Value = '1.794.07'
$ws.Cells.Item(16, 5).Value = '  +0.62%  '
$c = $ws.Cells.Item(17, 4); $c.NumberFormat = "@"; $c.Value = '90.31'
$ws.Cells.Item(17, 5).Value = '  +0.88%  '
$c = $ws.Cells.Item(18, 4); $c.NumberFormat = "@"; $c.Value = '0.00001064'
$ws.Cells.Item(18, 5).Value = '  +0.73%  '
$c = $ws.Cells.Item(19, 4); $c.NumberFormat = "@"; $c.Value = '0.06442'
$ws.Cells.Item(19, 5).Value = '  +0.14%  '
$c = $ws.Cells.Item(20, 4); $c.NumberFormat = "@"; $c.Value = '1.002'
$ws.Cells.Item(20, 5).Value = '  +0.23%  '
$c = $ws.Cells.Item(21, 4); $c.NumberFormat = "@"; $c.Value = '17.27'
$ws.Cells.Item(21, 5).Value = '  +2.87%  '
$c = $ws.Cells.Item(22, 4); $c.NumberFormat = "@"; $c.Value = '5.913'
$ws.Cells.Item(22, 5).Value = '  +0.11%  '
$ws.Cells.Item(23, 4).Value = '28.586.31'
$ws.Cells.Item(23, 5).Value = '  +4.29%  '
$c = $ws.Cells.Item(24, 4); $c.NumberFormat = "@"; $c.Value = '11.21'
$ws.Cells.Item(24, 5).Value = '  +0.13%  '
$c = $ws.Cells.Item(25, 4); $c.NumberFormat = "@"; $c.Value = '2.111'
$ws.Cells.Item(25, 5).Value = '  +1.05%  '
$c = $ws.Cells.Item(26, 4); $c.NumberFormat = "@"; $c.Value = '160.85'
$ws.Cells.Item(26, 5).Value = '  +3.43%  '
$c = $ws.Cells.Item(27, 4); $c.NumberFormat = "@"; $c.Value = '20.49'
$ws.Cells.Item(27, 5).Value = '  +1.39%  '
$c = $ws.Cells.Item(28, 4); $c.NumberFormat = "@"; $c.Value = '2.367'
$ws.Cells.Item(28, 5).Value = '  -0.30%  '
$ws.Cells.Item(29, 4).Value = '2.000.29'
$ws.Cells.Item(29, 5).Value = '  +0.65%  '
$c = $ws.Cells.Item(30, 4); $c.NumberFormat = "@"; $c.Value = '123.24'
$ws.Cells.Item(30, 5).Value = '  +1.75%  '
$c = $ws.Cells.Item(31, 4); $c.NumberFormat = "@"; $c.Value = '1.119'
$ws.Cells.Item(31, 5).Value = '  +3.84%  '
$ws.Cells.Item(32, 5).Value = '  -1.21%  '
$c = $ws.Cells.Item(33, 4); $c.NumberFormat = "@"; $c.Value = '5.693'
$ws.Cells.Item(33, 5).Value = '  +1.87%  '
$c = $ws.Cells.Item(34, 4); $c.NumberFormat = "@"; $c.Value = '3.659'
$ws.Cells.Item(34, 5).Value = '  +1.00%  '
$c = $ws.Cells.Item(35, 4); $c.NumberFormat = "@"; $c.Value = '0.2288'
$ws.Cells.Item(35, 5).Value = '  +11.21%  '
$c = $ws.Cells.Item(36, 4); $c.NumberFormat = "@"; $c.Value = '0.06561'
$ws.Cells.Item(36, 5).Value = '  +10.00%  '
$c = $ws.Cells.Item(37, 4); $c.NumberFormat = "@"; $c.Value = '0.02324'
$ws.Cells.Item(37, 5).Value = '  +2.82%  '
$ws.Cells.Item(38, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Cells.Item(38, 4); $c.NumberFormat = "@"; $c.Value = '5.088'
$ws.Cells.Item(38, 5).Value = '  +3.32%  '
$ws.Cells.Item(39, 2).Value = 'FraxShare'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(39, 4); $c.NumberFormat = "@"; $c.Value = '8.682'
$ws.Cells.Item(39, 5).Value = '  +5.20%  '
$c = $ws.Cells.Item(40, 4); $c.NumberFormat = "@"; $c.Value = '11.47'
$ws.Cells.Item(40, 5).Value = '  +1.85%  '
$c = $ws.Cells.Item(41, 4); $c.NumberFormat = "@"; $c.Value = '0.6314'
$ws.Cells.Item(41, 5).Value = '  +3.17%  '
$c = $ws.Cells.Item(42, 4); $c.NumberFormat = "@"; $c.Value = '1.204'
$ws.Cells.Item(42, 5).Value = '  +6.14%  '
$ws.Cells.Item(43, 5).Value = '  +0.32%  '
$c = $ws.Cells.Item(44, 4); $c.NumberFormat = "@"; $c.Value = '1.381'
$ws.Cells.Item(44, 5).Value = '  -3.37%  '
$c = $ws.Cells.Item(45, 4); $c.NumberFormat = "@"; $c.Value = '13.50'
$ws.Cells.Item(45, 5).Value = '  +1.30%  '
$c = $ws.Cells.Item(46, 4); $c.NumberFormat = "@"; $c.Value = '0.5917'
$ws.Cells.Item(46, 5).Value = '  +2.40%  '
$c = $ws.Cells.Item(47, 4); $c.NumberFormat = "@"; $c.Value = '3.667'
$ws.Cells.Item(47, 5).Value = '  +1.30%  '
$c = $ws.Cells.Item(48, 4); $c.NumberFormat = "@"; $c.Value = '125.43'
$ws.Cells.Item(48, 5).Value = '  +3.40%  '
$c = $ws.Cells.Item(49, 4); $c.NumberFormat = "@"; $c.Value = '1.978'
$ws.Cells.Item(49, 5).Value = '  +4.48%  '
$c = $ws.Cells.Item(50, 4); $c.NumberFormat = "@"; $c.Value = '1.163'
$ws.Cells.Item(50, 5).Value = '  +3.89%  '
$c = $ws.Cells.Item(51, 4); $c.NumberFormat = "@"; $c.Value = '0.06922'
$ws.Cells.Item(51, 5).Value = '  +2.95%  '
